# Venture Scheduler sample - update order row 6 (orderId, part, status label, quantity)
# and move the active selection to A6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 data updates
$ws.Range("A6").Value = 9874
$ws.Range("C6").Value = "test review"
$ws.Range("F6").Value = 1234

# Update the active cell/selection shown in the sheet view
$ws.Range("A6").Select()
